$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new "Narration" header in C1 (adds new shared string + fills
# previously-empty styled cell)
$ws.Range("C1").Value = "Narration"

# Reflect the active selection moving to C2, as captured in the saved file
$ws.Range("C2").Select()
